# Updated IPS AIP hipo turnover
$wb = $excel.ActiveWorkbook

# --- Sheet "Kongegårdsgatan Molndal Sweden": Internal Fill Rate / Commit-Forecast row (row 5) ---
$wsSweden = $wb.Worksheets.Item("Kongegårdsgatan Molndal Sweden")
$wsSweden.Range("M5").Value = $null   # Jun: 0.5 -> (blank)
$wsSweden.Range("N5").Value = 0       # Q2:  0.5 -> 0

# --- Sheet "Charlotte  North Carolina": Professional Voluntary Turnover / Commit-Forecast row (row 2) ---
$wsCharlotte = $wb.Worksheets.Item("Charlotte  North Carolina")
$wsCharlotte.Range("E2").Value = 0.6667               # ytd
$wsCharlotte.Range("M2").Value = 0                    # Jun
$wsCharlotte.Range("N2").Value = 1                    # Q2
$wsCharlotte.Range("O2").Value = 0.111116666666667    # Jul
$wsCharlotte.Range("P2").Value = 0.111116666666667    # Aug
$wsCharlotte.Range("Q2").Value = 0.111116666666667    # Sep
$wsCharlotte.Range("R2").Value = 0.33335              # Q3
$wsCharlotte.Range("S2").Value = 0.111116666666667    # Oct
$wsCharlotte.Range("T2").Value = 0.111116666666667    # Nov
$wsCharlotte.Range("U2").Value = 0.111116666666667    # Dec
$wsCharlotte.Range("V2").Value = 0.33335              # Q4
$wsCharlotte.Range("W2").Value = 1.3334               # FY

# --- Sheet "Shanghai Minhang District Chin": Internal Fill Rate / Commit-Forecast row (row 3) ---
$wsShanghai = $wb.Worksheets.Item("Shanghai Minhang District Chin")
$wsShanghai.Range("M3").Value = $null  # Jun: 1 -> (blank)
